$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 6.132086
$ws.Cells.Item(2,8).Value = 18.396258
$ws.Cells.Item(2,9).Value = 0.3910602616134352
$ws.Cells.Item(2,10).Value = 0.3910602616134352
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 6.132086
$ws.Cells.Item(2,14).Value = 18.396258
$ws.Cells.Item(2,15).Value = 0.3910602616134352
$ws.Cells.Item(2,16).Value = 0.3910602616134352
$ws.Cells.Item(2,17).Value = 37.602478711396
$ws.Cells.Item(2,18).Value = 338.422308402564
$ws.Cells.Item(2,19).Value = 0.1529281282131683
$ws.Cells.Item(2,20).Value = 0.1529281282131683
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 6.132086
$ws.Cells.Item(3,8).Value = 18.396258
$ws.Cells.Item(3,9).Value = 0.3910602616134352
$ws.Cells.Item(3,10).Value = 0.3910602616134352
$ws.Cells.Item(3,15).Value = 0.4122361243844145
$ws.Cells.Item(3,16).Value = 0.4122361243844145
$ws.Cells.Item(3,17).Value = 39.63864808783933
$ws.Cells.Item(3,18).Value = 356.747832790554
$ws.Cells.Item(3,19).Value = 0.1612091666482777
$ws.Cells.Item(3,20).Value = 0.1612091666482777
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 6.132086
$ws.Cells.Item(4,8).Value = 18.396258
$ws.Cells.Item(4,9).Value = 0.3910602616134352
$ws.Cells.Item(4,10).Value = 0.3910602616134352
$ws.Cells.Item(4,13).Value = 2.936401666666667
$ws.Cells.Item(4,14).Value = 8.809205
$ws.Cells.Item(4,15).Value = 0.1872625406703027
$ws.Cells.Item(4,16).Value = 0.1872625406703027
$ws.Cells.Item(4,17).Value = 18.00626755054333
$ws.Cells.Item(4,18).Value = 162.05640795489
$ws.Cells.Item(4,19).Value = 0.07323093814492511
$ws.Cells.Item(4,20).Value = 0.07323093814492512
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 6.132086
$ws.Cells.Item(5,8).Value = 18.396258
$ws.Cells.Item(5,9).Value = 0.3910602616134352
$ws.Cells.Item(5,10).Value = 0.3910602616134352
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.1480423333333333
$ws.Cells.Item(5,14).Value = 0.444127
$ws.Cells.Item(5,15).Value = 0.00944107333184771
$ws.Cells.Item(5,16).Value = 0.009441073331847712
$ws.Cells.Item(5,17).Value = 0.9078083196406667
$ws.Cells.Item(5,18).Value = 8.170274876765999
$ws.Cells.Item(5,19).Value = 0.003692028607063992
$ws.Cells.Item(5,20).Value = 0.003692028607063992
$ws.Cells.Item(6,9).Value = 0.4122361243844145
$ws.Cells.Item(6,10).Value = 0.4122361243844145
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 6.132086
$ws.Cells.Item(6,14).Value = 18.396258
$ws.Cells.Item(6,15).Value = 0.3910602616134352
$ws.Cells.Item(6,16).Value = 0.3910602616134352
$ws.Cells.Item(6,17).Value = 39.63864808783933
$ws.Cells.Item(6,18).Value = 356.747832790554
$ws.Cells.Item(6,19).Value = 0.1612091666482777
$ws.Cells.Item(6,20).Value = 0.1612091666482777
$ws.Cells.Item(7,9).Value = 0.4122361243844145
$ws.Cells.Item(7,10).Value = 0.4122361243844145
$ws.Cells.Item(7,15).Value = 0.4122361243844145
$ws.Cells.Item(7,16).Value = 0.4122361243844145
$ws.Cells.Item(7,19).Value = 0.1699386222474824
$ws.Cells.Item(7,20).Value = 0.1699386222474824
$ws.Cells.Item(8,9).Value = 0.4122361243844145
$ws.Cells.Item(8,10).Value = 0.4122361243844145
$ws.Cells.Item(8,13).Value = 2.936401666666667
$ws.Cells.Item(8,14).Value = 8.809205
$ws.Cells.Item(8,15).Value = 0.1872625406703027
$ws.Cells.Item(8,16).Value = 0.1872625406703027
$ws.Cells.Item(8,17).Value = 18.98130461796278
$ws.Cells.Item(8,18).Value = 170.831741561665
$ws.Cells.Item(8,19).Value = 0.07719638400830438
$ws.Cells.Item(8,20).Value = 0.07719638400830438
$ws.Cells.Item(9,9).Value = 0.4122361243844145
$ws.Cells.Item(9,10).Value = 0.4122361243844145
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.1480423333333333
$ws.Cells.Item(9,14).Value = 0.444127
$ws.Cells.Item(9,15).Value = 0.00944107333184771
$ws.Cells.Item(9,16).Value = 0.009441073331847712
$ws.Cells.Item(9,17).Value = 0.9569660231612221
$ws.Cells.Item(9,18).Value = 8.612694208450998
$ws.Cells.Item(9,19).Value = 0.003891951480349951
$ws.Cells.Item(9,20).Value = 0.003891951480349952
$ws.Cells.Item(10,7).Value = 2.936401666666667
$ws.Cells.Item(10,8).Value = 8.809205
$ws.Cells.Item(10,9).Value = 0.1872625406703027
$ws.Cells.Item(10,10).Value = 0.1872625406703027
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 6.132086
$ws.Cells.Item(10,14).Value = 18.396258
$ws.Cells.Item(10,15).Value = 0.3910602616134352
$ws.Cells.Item(10,16).Value = 0.3910602616134352
$ws.Cells.Item(10,17).Value = 18.00626755054333
$ws.Cells.Item(10,18).Value = 162.05640795489
$ws.Cells.Item(10,19).Value = 0.07323093814492511
$ws.Cells.Item(10,20).Value = 0.07323093814492512
$ws.Cells.Item(11,7).Value = 2.936401666666667
$ws.Cells.Item(11,8).Value = 8.809205
$ws.Cells.Item(11,9).Value = 0.1872625406703027
$ws.Cells.Item(11,10).Value = 0.1872625406703027
$ws.Cells.Item(11,15).Value = 0.4122361243844145
$ws.Cells.Item(11,16).Value = 0.4122361243844145
$ws.Cells.Item(11,17).Value = 18.98130461796278
$ws.Cells.Item(11,18).Value = 170.831741561665
$ws.Cells.Item(11,19).Value = 0.07719638400830438
$ws.Cells.Item(11,20).Value = 0.07719638400830438
$ws.Cells.Item(12,7).Value = 2.936401666666667
$ws.Cells.Item(12,8).Value = 8.809205
$ws.Cells.Item(12,9).Value = 0.1872625406703027
$ws.Cells.Item(12,10).Value = 0.1872625406703027
$ws.Cells.Item(12,13).Value = 2.936401666666667
$ws.Cells.Item(12,14).Value = 8.809205
$ws.Cells.Item(12,15).Value = 0.1872625406703027
$ws.Cells.Item(12,16).Value = 0.1872625406703027
$ws.Cells.Item(12,17).Value = 8.622454748002779
$ws.Cells.Item(12,18).Value = 77.60209273202501
$ws.Cells.Item(12,19).Value = 0.03506725913829677
$ws.Cells.Item(12,20).Value = 0.03506725913829677
$ws.Cells.Item(13,7).Value = 2.936401666666667
$ws.Cells.Item(13,8).Value = 8.809205
$ws.Cells.Item(13,9).Value = 0.1872625406703027
$ws.Cells.Item(13,10).Value = 0.1872625406703027
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 0.1480423333333333
$ws.Cells.Item(13,14).Value = 0.444127
$ws.Cells.Item(13,15).Value = 0.00944107333184771
$ws.Cells.Item(13,16).Value = 0.009441073331847712
$ws.Cells.Item(13,17).Value = 0.4347117543372223
$ws.Cells.Item(13,18).Value = 3.912405789035
$ws.Cells.Item(13,19).Value = 0.001767959378776442
$ws.Cells.Item(13,20).Value = 0.001767959378776442
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 0.1480423333333333
$ws.Cells.Item(14,8).Value = 0.444127
$ws.Cells.Item(14,9).Value = 0.00944107333184771
$ws.Cells.Item(14,10).Value = 0.009441073331847712
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 6.132086
$ws.Cells.Item(14,14).Value = 18.396258
$ws.Cells.Item(14,15).Value = 0.3910602616134352
$ws.Cells.Item(14,16).Value = 0.3910602616134352
$ws.Cells.Item(14,17).Value = 0.9078083196406667
$ws.Cells.Item(14,18).Value = 8.170274876765999
$ws.Cells.Item(14,19).Value = 0.003692028607063992
$ws.Cells.Item(14,20).Value = 0.003692028607063992
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 0.1480423333333333
$ws.Cells.Item(15,8).Value = 0.444127
$ws.Cells.Item(15,9).Value = 0.00944107333184771
$ws.Cells.Item(15,10).Value = 0.009441073331847712
$ws.Cells.Item(15,15).Value = 0.4122361243844145
$ws.Cells.Item(15,16).Value = 0.4122361243844145
$ws.Cells.Item(15,17).Value = 0.9569660231612221
$ws.Cells.Item(15,18).Value = 8.612694208450998
$ws.Cells.Item(15,19).Value = 0.003891951480349951
$ws.Cells.Item(15,20).Value = 0.003891951480349952
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 0.1480423333333333
$ws.Cells.Item(16,8).Value = 0.444127
$ws.Cells.Item(16,9).Value = 0.00944107333184771
$ws.Cells.Item(16,10).Value = 0.009441073331847712
$ws.Cells.Item(16,13).Value = 2.936401666666667
$ws.Cells.Item(16,14).Value = 8.809205
$ws.Cells.Item(16,15).Value = 0.1872625406703027
$ws.Cells.Item(16,16).Value = 0.1872625406703027
$ws.Cells.Item(16,17).Value = 0.4347117543372223
$ws.Cells.Item(16,18).Value = 3.912405789035
$ws.Cells.Item(16,19).Value = 0.001767959378776442
$ws.Cells.Item(16,20).Value = 0.001767959378776442
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 0.1480423333333333
$ws.Cells.Item(17,8).Value = 0.444127
$ws.Cells.Item(17,9).Value = 0.00944107333184771
$ws.Cells.Item(17,10).Value = 0.009441073331847712
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 0.1480423333333333
$ws.Cells.Item(17,14).Value = 0.444127
$ws.Cells.Item(17,15).Value = 0.00944107333184771
$ws.Cells.Item(17,16).Value = 0.009441073331847712
$ws.Cells.Item(17,17).Value = 0.02191653245877778
$ws.Cells.Item(17,18).Value = 0.197248792129
$ws.Cells.Item(17,19).Value = 0.00008913386565732603
$ws.Cells.Item(17,20).Value = 0.00008913386565732606
